$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value2
$text = $text.Replace("1000 Bs = 6.99 = 27378.58 pesos", "1000 Bs = 6.98 = 27308.1 pesos")
$text = $text.Replace("27378.58 pesos = 6.97 = 974.68 Bs", "27308.1 pesos = 6.95 = 966.55 Bs")
$cell.Value = $text

# --- Sheet "tasas": update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 143.2
$wsTasas.Range("O10").Value = 3910.52
$wsTasas.Range("N12").Value = 3930
$wsTasas.Range("O12").Value = 139.1
